$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "DP on Matrix" column (E) with its two entries
$ws.Range("E1").Value2 = "DP on Matrix"
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").Font.Size = 12

$ws.Range("E2").Value2 = "Leetcode - 200"
$ws.Range("E3").Value2 = "Leetcode - 221"

# Update the existing header in column A: "Maximum SubArray" -> "Maximum/Minimum SubArray"
$ws.Range("A1").Value2 = "Maximum/Minimum SubArray"

# Size column E similarly to the other bestFit columns
$ws.Columns("E:E").ColumnWidth = 12.3

# Move the active selection to C16 (matches the saved workbook view)
$ws.Range("C16").Select() | Out-Null
